# Update to v1.1.0: rework "Fund Transaction Report" (dividends) sheet columns.
#  - remove the "分红交易所属日" / {.dividendDate} column
#  - add "分红次数" / {.dividendCount} and "合计分红金额" / {.totalDividendAmount} columns
#    (placed right after 合计金额 / {.totalAmount})
#  - add "收益率" / {.yieldRate} column (placed right after 日均万份收益 / {.dailyNavYield})

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fund Transaction Report")

# 1. Delete the "分红交易所属日" column (old column E).
$ws.Columns.Item(5).Delete()

# After the delete the columns are:
#   A 编号 | B 基金简称 | C 购买交易所属日 | D 赎回交易所属日 | E 持有天数 |
#   F 合计本金 | G 合计金额 | H 净收益 | I 日均万份收益 | J 交易平台 | K 基金公司

# 2. Insert two new columns right before the old "净收益" column (H) for the
#    dividend-count / total-dividend-amount fields.
$ws.Range("H1:I1").EntireColumn.Insert()

# Columns are now:
#   A..G (unchanged) | H (new) | I (new) | J 净收益 | K 日均万份收益 | L 交易平台 | M 基金公司

# 3. Insert one new column right before "交易平台" (now column L) for the yield-rate field.
$ws.Range("L1").EntireColumn.Insert()

# Final layout:
#   A 编号 | B 基金简称 | C 购买交易所属日 | D 赎回交易所属日 | E 持有天数 | F 合计本金 |
#   G 合计金额 | H 分红次数 | I 合计分红金额 | J 净收益 | K 日均万份收益 | L 收益率 |
#   M 交易平台 | N 基金公司

# 4. Populate the header (row 2) / template-placeholder (row 3) text for the new columns.
$ws.Range("H2").Value = "分红次数"
$ws.Range("H3").Value = "{.dividendCount}"
$ws.Range("L2").Value = "收益率"
$ws.Range("L3").Value = "{.yieldRate}"
$ws.Range("I3").Value = "{.totalDividendAmount}"
$ws.Range("I2").Value = "合计分红金额"

# 5. Match the target column widths as closely as the Excel column-width
#    pixel-grid allows.
$ws.Columns.Item(5).ColumnWidth = 10.1          # heldDays column -> width 11
$ws.Columns.Item(9).ColumnWidth = 13.8          # totalDividendAmount column -> width ~14.71
$ws.Columns.Item(12).ColumnWidth = 10.1         # yieldRate column -> width 11
